$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.02   # Current Capital
$summary.Range("B4").Value = 0.02      # Total P&L $
$summary.Range("B5").Value = 0.4       # Total P&L %
$summary.Range("B6").Value = 1         # Total Trades
$summary.Range("B7").Value = 1         # Winning Trades
$summary.Range("B9").Value = 100       # Win Rate %

# --- Strategy Status sheet (MarketMaking row, row 4) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.02     # Capital
$status.Range("D4").Value = 1          # Trades
$status.Range("E4").Value = 0.02       # P&L $
$status.Range("F4").Value = 0.02       # P&L %
$status.Range("G4").Value = 100        # Win Rate %

# --- All Trades & MarketMaking sheets share the same trade-row update ---
$sheetsToUpdate = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetsToUpdate) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2").Value = "07:57:02"
    $ws.Range("F2").Value = 0.15
    $ws.Range("G2").Value = 0.17
    $ws.Range("H2").Value = "CLOSED"
    $ws.Range("I2").Value = 13.3333
    $ws.Range("J2").Value = 0.02
    $ws.Range("K2").Value = 100.02
    $ws.Range("P2").Value = "early_exit"
    $ws.Range("Q2").Value = 0.11
}
